# liensMagasinPT.xlsx - "IT" store links renamed to "PT" store links.
# Column A (rows 1-44) holds a repeated section-title string per block of
# rows; every occurrence of the old "...ClientIT" / "...FullWebIT" /
# "Magasin Citrix ClientIT" title is replaced by its "PT" counterpart
# (the "Citrix" block also gets reworded to "CPTrix").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A4").Value = "Magasin CL Kstore ClientPT"
$ws.Range("A5:A13").Value = "Magasin CLO ClientPT"
$ws.Range("A14:A20").Value = "Magasin FullWebPT"
$ws.Range("A21").Value = "Magasin LSA ClientPT"
$ws.Range("A22:A33").Value = "Magasin POD ClientPT"
$ws.Range("A34:A40").Value = "Magasin SCO ClientPT"
$ws.Range("A41:A44").Value = "Magasin CPTrix ClientPT"

# Restore the view state as closely as the object model allows: the
# sheet had scrolled so row 33 is at the top, with A1:A44 selected and
# A44 as the active cell.
$ws.Range("A1:A44").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
